$wb = $excel.ActiveWorkbook

# Update the Users sheet: replace "Gemma Hardy" with the new user "Jennie Stewart"
$wsUsers = $wb.Worksheets.Item("Users")
$wsUsers.Range("B2").Value = "Jennie Stewart"

# Make Users the active sheet/tab with B2 selected (mirrors the saved view state)
$wsUsers.Activate() | Out-Null
$wsUsers.Range("B2").Select() | Out-Null
